# Automatische test-sync: 2025-07-27 19:21:50
# Append the 5th test-mail log entry to the "Logs" sheet and bump the
# matching category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append new row 7 to the Logs sheet -----------------------------------
$newRow = 7

$logs.Cells.Item($newRow, 1).Value = "Wil je deze klant bellen?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #5: Wil je deze klant bellen?"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value = "Geachte heer/mevrouw,`nDank u voor uw bericht. Helaas kan ik niet achterhalen om welke klant het gaat op basis van de informatie die u heeft verstrekt. Kunt u mij de naam of het klantnummer van de desbetreffende klant geven, zodat ik dit verder kan onderzoeken?`nMet vriendelijke groet,`n[Naam van de e-mailassistent]`n[Naam van het bedrijf]"
$logs.Cells.Item($newRow, 6).Value = "2025-07-27 19:21:11"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Writing the multi-line comment above causes the host to auto-grow the new
# row's height; AutoFit() puts it back to the sheet's standard height so the
# row doesn't end up with a stray explicit/custom height like a real user
# typing the same content wouldn't produce either.
$logs.Rows.Item($newRow).AutoFit()

# --- Extend conditional formatting ranges to cover the new row ------------
function Extend-ConditionalFormatting($columnLetter, $oldLastRow, $newLastRow) {
    $oldRange = $logs.Range($columnLetter + "2:" + $columnLetter + $oldLastRow)
    $newRange = $logs.Range($columnLetter + "2:" + $columnLetter + $newLastRow)
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

Extend-ConditionalFormatting "D" 6 7
Extend-ConditionalFormatting "G" 6 7
Extend-ConditionalFormatting "H" 6 7
Extend-ConditionalFormatting "I" 6 7
Extend-ConditionalFormatting "J" 6 7

# --- Update the Dashboard summary count for the matching category ---------
$dashboard.Cells.Item(3, 2).Value = 2
